# Opencart_LoginData.xlsx - refresh the test credentials used by the
# login test (old/invalid-login test account is being swapped for a new
# one, and the WebDriver "status" bookkeeping that used to live here is
# no longer needed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old invalid-login test account (row 2) with the new one.
$ws.Range("A2").Value = "tejaszombade55@gmail.com"
$ws.Range("B2").Value = "Tztejas@13"

# The new values pick up the compact "hyperlink" look used by the other
# rows in the sheet (11pt instead of 16pt), which also shrinks the row.
$ws.Range("A2:B2").Font.Size = 11
$ws.Rows.Item(2).RowHeight = 21

# Leave the cursor where it was left after the edit.
$ws.Range("E6").Select() | Out-Null
